$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the first data row (row 2), pushing the existing
# fiscal-year rows (2020, 2021, 2022) down by one.
$ws.Rows("2:2").Insert()

# The surrounding data rows store name/date_from/date_to as text (e.g.
# "2020", "2020-01-01"), so force those new cells to text format before
# writing the values to avoid Excel auto-converting "2019" to a number or
# the dates to date serials. Column A (id) and E (status) keep the default
# general style, matching the other data rows.
$ws.Range("B2:D2").NumberFormat = "@"

# Fill in the new fiscal year 2019 row.
$ws.Range("A2").Value = "z0bug.fy_2019"
$ws.Range("B2").Value = "2019"
$ws.Range("C2").Value = "2019-01-01"
$ws.Range("D2").Value = "2019-12-31"
$ws.Range("E2").Value = "Open"

$ws.Range("D3").Select()
